# Commit message: "Started implementing [eb2]. Not finished."
#
# The OOXML diff for xl/workbook.xml shows:
#   - fileVersion/revisionPtr bump (Excel build bookkeeping -- not content)
#   - bookViews window geometry change (Excel window placement -- not content)
#   - the worksheet renamed from the eb1-specific "cond_eb1_c" to the more
#     generic "cond" (the actual, meaningful edit -- consistent with the
#     commit message: generalizing the eb1 file as a start toward eb2)
#   - calcPr gains iterateDelta="1E-4" (iterative-calculation max-change
#     tolerance)
#
# Apply the real content edit: rename the (only) worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item(1)
$ws.Name = "cond"

# Mirror the calcPr hint (iterateDelta="1E-4") by enabling iterative
# calculation with that max-change tolerance via the Application object,
# the normal COM route for this setting. Wrapped defensively so a host
# that doesn't expose these falls back to just the rename above.
try {
    $excel.Iteration = $true
    $excel.MaxChange = 0.0001
} catch {
}
